# Worked on temporal resolution:
# Split the single annual "Demand" value for EU27.Elec (t=1) into 12 equal
# per-period values (t=1..12), matching a finer temporal resolution.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Annual total (previously stored at B3) split evenly across 12 periods.
$periodValue = 406359375

for ($t = 1; $t -le 12; $t++) {
    $row = $t + 2  # row 3 holds t=1, row 4 holds t=2, ... row 14 holds t=12
    $ws.Cells.Item($row, 1).Value = $t
    $ws.Cells.Item($row, 2).Value = $periodValue
}

# Column B now holds a wider / longer series of numbers - widen it to fit.
$ws.Columns.Item(2).ColumnWidth = 9.15

# Bring the Demand sheet to the front / make it the active tab & selection.
$ws.Activate() | Out-Null
$ws.Range("I12").Select() | Out-Null
